$d = $word.ActiveDocument

# 1) Hearing date: April 06, 2022 -> April 04, 2022
$d.Content.Find.Execute("April 06, 2022", $true, $true, $false, $false, $false,
                         $true, 1, $false, "April 04, 2022", 2) | Out-Null

# 2) Hearing purpose: "for change of plea" -> "for arraignment"
$d.Content.Find.Execute(" for change of plea", $true, $true, $false, $false, $false,
                         $true, 1, $false, " for arraignment", 2) | Out-Null

# 3) Statute/Ord. row: 4511.19A2*** -> 4511.19A1A***
$d.Content.Find.Execute("4511.19A2***", $true, $true, $false, $false, $false,
                         $true, 1, $false, "4511.19A1A***", 2) | Out-Null

# 4) Statute/Ord. row: 4511.33 -> 4511.19A1A***
$d.Content.Find.Execute("4511.33", $true, $true, $false, $false, $false,
                         $true, 1, $false, "4511.19A1A***", 2) | Out-Null

# 5) Statute/Ord. row: 4513.263B1 -> 4511.19A1A***
$d.Content.Find.Execute("4513.263B1", $true, $true, $false, $false, $false,
                         $true, 1, $false, "4511.19A1A***", 2) | Out-Null

# 6) Degree row: MM -> UCM
$d.Content.Find.Execute("MM", $true, $true, $true, $false, $false,
                         $true, 1, $false, "UCM", 2) | Out-Null

# 7) Add keepNext to the paragraph holding the signature underscore line,
#    directly above "Judge ... " so it stays with the following paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "____________________________________") {
        $p.Range.ParagraphFormat.KeepWithNext = $true
        break
    }
}

# 8) Judge name: Kyle -> Marianne
$d.Content.Find.Execute("Kyle", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Marianne", 2) | Out-Null

# 9) Judge surname: Rohrer -> Hemmeter
$d.Content.Find.Execute("Rohrer", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Hemmeter", 2) | Out-Null
